$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 12.31940633333333
$ws.Range("H2").Value = 36.958219
$ws.Range("I2").Value = 0.3801768364207869
$ws.Range("J2").Value = 0.3801768364207869
$ws.Range("M2").Value = 91.60947133333333
$ws.Range("N2").Value = 274.828414
$ws.Range("O2").Value = 0.7121576185577153
$ws.Range("P2").Value = 0.7121576185577152
$ws.Range("Q2").Value = 1128.574301337185
$ws.Range("R2").Value = 10157.16871203467
$ws.Range("S2").Value = 0.2707458304562337
$ws.Range("T2").Value = 0.2707458304562336
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 12.31940633333333
$ws.Range("H3").Value = 36.958219
$ws.Range("I3").Value = 0.3801768364207869
$ws.Range("J3").Value = 0.3801768364207869
$ws.Range("M3").Value = 22.83185066666667
$ws.Range("N3").Value = 68.495552
$ws.Range("O3").Value = 0.1774912152792038
$ws.Range("P3").Value = 0.1774912152792038
$ws.Range("Q3").Value = 281.2748457046542
$ws.Range("R3").Value = 2531.473611341888
$ws.Range("S3").Value = 0.06747804871732854
$ws.Range("T3").Value = 0.06747804871732853
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.31940633333333
$ws.Range("H4").Value = 36.958219
$ws.Range("I4").Value = 0.3801768364207869
$ws.Range("J4").Value = 0.3801768364207869
$ws.Range("M4").Value = 7.077809999999999
$ws.Range("N4").Value = 21.23343
$ws.Range("O4").Value = 0.05502178149094856
$ws.Range("P4").Value = 0.05502178149094855
$ws.Range("Q4").Value = 87.19441734012999
$ws.Range("R4").Value = 784.7497560611699
$ws.Range("S4").Value = 0.02091800682146463
$ws.Range("T4").Value = 0.02091800682146462
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.31940633333333
$ws.Range("H5").Value = 36.958219
$ws.Range("I5").Value = 0.3801768364207869
$ws.Range("J5").Value = 0.3801768364207869
$ws.Range("M5").Value = 7.117379
$ws.Range("N5").Value = 21.352137
$ws.Range("O5").Value = 0.05532938467213248
$ws.Range("P5").Value = 0.05532938467213247
$ws.Range("Q5").Value = 87.68188392933367
$ws.Range("R5").Value = 789.1369553640029
$ws.Range("S5").Value = 0.0210349504257601
$ws.Range("T5").Value = 0.0210349504257601
$ws.Range("G6").Value = 18.11265066666666
$ws.Range("H6").Value = 54.33795199999999
$ws.Range("I6").Value = 0.5589563363143816
$ws.Range("J6").Value = 0.5589563363143816
$ws.Range("M6").Value = 91.60947133333333
$ws.Range("N6").Value = 274.828414
$ws.Range("O6").Value = 0.7121576185577153
$ws.Range("P6").Value = 0.7121576185577152
$ws.Range("Q6").Value = 1659.290352018681
$ws.Range("R6").Value = 14933.61316816813
$ws.Range("S6").Value = 0.3980650133473954
$ws.Range("T6").Value = 0.3980650133473953
$ws.Range("G7").Value = 18.11265066666666
$ws.Range("H7").Value = 54.33795199999999
$ws.Range("I7").Value = 0.5589563363143816
$ws.Range("J7").Value = 0.5589563363143816
$ws.Range("M7").Value = 22.83185066666667
$ws.Range("N7").Value = 68.495552
$ws.Range("O7").Value = 0.1774912152792038
$ws.Range("P7").Value = 0.1774912152792038
$ws.Range("S7").Value = 0.09920983942045096
$ws.Range("T7").Value = 0.09920983942045095
$ws.Range("G8").Value = 18.11265066666666
$ws.Range("H8").Value = 54.33795199999999
$ws.Range("I8").Value = 0.5589563363143816
$ws.Range("J8").Value = 0.5589563363143816
$ws.Range("M8").Value = 7.077809999999999
$ws.Range("N8").Value = 21.23343
$ws.Range("O8").Value = 0.05502178149094856
$ws.Range("P8").Value = 0.05502178149094855
$ws.Range("Q8").Value = 128.19790001504
$ws.Range("R8").Value = 1153.78110013536
$ws.Range("S8").Value = 0.03075477339967106
$ws.Range("T8").Value = 0.03075477339967105
$ws.Range("G9").Value = 18.11265066666666
$ws.Range("H9").Value = 54.33795199999999
$ws.Range("I9").Value = 0.5589563363143816
$ws.Range("J9").Value = 0.5589563363143816
$ws.Range("M9").Value = 7.117379
$ws.Range("N9").Value = 21.352137
$ws.Range("O9").Value = 0.05532938467213248
$ws.Range("P9").Value = 0.05532938467213247
$ws.Range("Q9").Value = 128.9145994892693
$ws.Range("R9").Value = 1160.231395403424
$ws.Range("S9").Value = 0.03092671014686427
$ws.Range("T9").Value = 0.03092671014686427
$ws.Range("G10").Value = 1.603212
$ws.Range("H10").Value = 4.809636
$ws.Range("I10").Value = 0.04947511672073613
$ws.Range("J10").Value = 0.04947511672073613
$ws.Range("M10").Value = 91.60947133333333
$ws.Range("N10").Value = 274.828414
$ws.Range("O10").Value = 0.7121576185577153
$ws.Range("P10").Value = 0.7121576185577152
$ws.Range("Q10").Value = 146.869403755256
$ws.Range("R10").Value = 1321.824633797304
$ws.Range("S10").Value = 0.03523408130170445
$ws.Range("T10").Value = 0.03523408130170444
$ws.Range("G11").Value = 1.603212
$ws.Range("H11").Value = 4.809636
$ws.Range("I11").Value = 0.04947511672073613
$ws.Range("J11").Value = 0.04947511672073613
$ws.Range("M11").Value = 22.83185066666667
$ws.Range("N11").Value = 68.495552
$ws.Range("O11").Value = 0.1774912152792038
$ws.Range("P11").Value = 0.1774912152792038
$ws.Range("Q11").Value = 36.604296971008
$ws.Range("R11").Value = 329.438672739072
$ws.Range("S11").Value = 0.008781398592843914
$ws.Range("T11").Value = 0.008781398592843912
$ws.Range("G12").Value = 1.603212
$ws.Range("H12").Value = 4.809636
$ws.Range("I12").Value = 0.04947511672073613
$ws.Range("J12").Value = 0.04947511672073613
$ws.Range("M12").Value = 7.077809999999999
$ws.Range("N12").Value = 21.23343
$ws.Range("O12").Value = 0.05502178149094856
$ws.Range("P12").Value = 0.05502178149094855
$ws.Range("Q12").Value = 11.34722992572
$ws.Range("R12").Value = 102.12506933148
$ws.Range("S12").Value = 0.002722209061447519
$ws.Range("T12").Value = 0.002722209061447518
$ws.Range("G13").Value = 1.603212
$ws.Range("H13").Value = 4.809636
$ws.Range("I13").Value = 0.04947511672073613
$ws.Range("J13").Value = 0.04947511672073613
$ws.Range("M13").Value = 7.117379
$ws.Range("N13").Value = 21.352137
$ws.Range("O13").Value = 0.05532938467213248
$ws.Range("P13").Value = 0.05532938467213247
$ws.Range("Q13").Value = 11.410667421348
$ws.Range("R13").Value = 102.696006792132
$ws.Range("S13").Value = 0.002737427764740263
$ws.Range("T13").Value = 0.002737427764740263
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.3691416666666666
$ws.Range("H14").Value = 1.107425
$ws.Range("I14").Value = 0.01139171054409548
$ws.Range("J14").Value = 0.01139171054409548
$ws.Range("M14").Value = 91.60947133333333
$ws.Range("N14").Value = 274.828414
$ws.Range("O14").Value = 0.7121576185577153
$ws.Range("P14").Value = 0.7121576185577152
$ws.Range("Q14").Value = 33.81687293043888
$ws.Range("R14").Value = 304.35185637395
$ws.Range("S14").Value = 0.008112693452381851
$ws.Range("T14").Value = 0.00811269345238185
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.3691416666666666
$ws.Range("H15").Value = 1.107425
$ws.Range("I15").Value = 0.01139171054409548
$ws.Range("J15").Value = 0.01139171054409548
$ws.Range("M15").Value = 22.83185066666667
$ws.Range("N15").Value = 68.495552
$ws.Range("O15").Value = 0.1774912152792038
$ws.Range("P15").Value = 0.1774912152792038
$ws.Range("Q15").Value = 8.428187408177777
$ws.Range("R15").Value = 75.8536866736
$ws.Range("S15").Value = 0.002021928548580427
$ws.Range("T15").Value = 0.002021928548580426
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.3691416666666666
$ws.Range("H16").Value = 1.107425
$ws.Range("I16").Value = 0.01139171054409548
$ws.Range("J16").Value = 0.01139171054409548
$ws.Range("M16").Value = 7.077809999999999
$ws.Range("N16").Value = 21.23343
$ws.Range("O16").Value = 0.05502178149094856
$ws.Range("P16").Value = 0.05502178149094855
$ws.Range("Q16").Value = 2.61271457975
$ws.Range("R16").Value = 23.51443121774999
$ws.Range("S16").Value = 0.0006267922083653561
$ws.Range("T16").Value = 0.000626792208365356
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.3691416666666666
$ws.Range("H17").Value = 1.107425
$ws.Range("I17").Value = 0.01139171054409548
$ws.Range("J17").Value = 0.01139171054409548
$ws.Range("M17").Value = 7.117379
$ws.Range("N17").Value = 21.352137
$ws.Range("O17").Value = 0.05532938467213248
$ws.Range("P17").Value = 0.05532938467213247
$ws.Range("Q17").Value = 2.627321146358333
$ws.Range("R17").Value = 23.645890317225
$ws.Range("S17").Value = 0.0006302963347678464
$ws.Range("T17").Value = 0.0006302963347678462
